$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 becomes the "Unassigned" row (was row 40, lowercase "unassigned")
$ws.Range("A39").Value = "Unassigned"
$ws.Range("B39").Value = "Unassigned"
$ws.Range("C39").Value = "Unassigned"
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 298
$ws.Range("F39").Value = 239

# Row 40 becomes the "Urophycis sp" row (was row 39)
$ws.Range("A40").Value = "Urophycis sp"
$ws.Range("B40").Value = "Red White or Spotted hake"
$ws.Range("C40").Value = "Teleost Fish"
$ws.Range("D40").Value = 1261
$ws.Range("E40").Value = 1025
$ws.Range("F40").Value = 11
